# Generate Report for Handoff
# Replaces the source-file GUID with a freshly generated one, records the
# new handoff (xlf) package hashes/timestamps, and resets the "handback"
# columns for each locale now that a brand-new handoff round has begun.

$wb = $excel.ActiveWorkbook

$oldGuid = "ad524f1c-9d01-4152-ac78-f7db533094fb"
$newGuid = "928eb79a-370c-46aa-800b-3ade8f28e38e"

$newFileName      = "$newGuid.md"
$newPathAndName   = "e2e\$newGuid.md"
$newHash          = "169c931c11bcad59b977493b44ea1d2f87f78249"
$zhHandoffXlf     = "$newGuid.$newHash.zh-cn.xlf"
$deHandoffXlf     = "$newGuid.$newHash.de-de.xlf"
$zhHandoffTime    = "2016-08-26 06:58:17"
$deHandoffTime    = "2016-08-26 06:58:22"
$resetHandback    = "0001-01-01 00:00:00"
$latestGenerate   = "2016-08-26 06:58:22"

function Clear-HandbackColumns($ws) {
    # Remove the hyperlink that lived on the "Latest Target File" cell (I2)
    # before touching the cell contents, mirroring what Excel does when a
    # handed-back file no longer exists.
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$I$2') {
            $hl.Delete()
        }
    }
    $ws.Range("I2").Style = "Normal"
    $ws.Range("I2").Value = ""
    $ws.Range("J2").Value = ""
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathAndName
foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.TextToDisplay = $newPathAndName
    }
}
$wsOverview.Range("G2").Value = $latestGenerate

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newFileName
foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = $newFileName
    }
}
$wsZh.Range("G2").Value = $zhHandoffXlf
$wsZh.Range("H2").Value = $zhHandoffTime
Clear-HandbackColumns $wsZh
$wsZh.Range("K2").Value = $resetHandback
$wsZh.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsZh.Columns.Item(10).ColumnWidth = 20.833333333333332

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newFileName
foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = $newFileName
    }
}
$wsDe.Range("G2").Value = $deHandoffXlf
$wsDe.Range("H2").Value = $deHandoffTime
Clear-HandbackColumns $wsDe
$wsDe.Range("K2").Value = $resetHandback
$wsDe.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsDe.Columns.Item(10).ColumnWidth = 20.833333333333332
